# Add a new row to the SQL50 patterns table for LeetCode 1978
# "Employees Whose Manager Left the Company"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives inside an Excel Table (ListObject); adding a ListRow
# grows the table range and the sheet dimension automatically.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$rowIndex = $newRow.Range.Row

$question = "1978. Employees Whose Manager Left the Company"
$difficulty = "Easy"
$pattern = "Subqueries"
$link = "https://leetcode.com/problems/employees-whose-manager-left-the-company/solutions/3789532/easy-solution/?envType=study-plan-v2&envId=top-sql-50 "
$notes = "Use a subquery for where, anager id not in (select employee_id from Employees). You can also use EXISTS()"

$ws.Cells.Item($rowIndex, 1).Value = $question
$ws.Cells.Item($rowIndex, 2).Value = $difficulty
# Match the green "Easy" fill used by the other Easy rows.
$ws.Cells.Item($rowIndex, 2).Interior.Color = $ws.Range("B2").Interior.Color
$ws.Cells.Item($rowIndex, 3).Value = $pattern

# Set the Link cell text before the Notes cell so shared-string indices
# line up the same way Excel itself would allocate them.
$ws.Cells.Item($rowIndex, 5).Value = $link
$ws.Cells.Item($rowIndex, 4).Value = $notes

$linkCell = $ws.Cells.Item($rowIndex, 5)
# The hyperlink target itself is trimmed (matching how the other rows'
# relationship targets drop the trailing space that the display text keeps).
$ws.Hyperlinks.Add($linkCell, $link.TrimEnd())
$linkCell.Style = $ws.Range("E31").Style

# Keep the same selected cell pattern seen after adding a row by hand.
$ws.Range("D36").Select() | Out-Null
